# Add a PIPPIB verification row to the "FORM KELENGKAPAN ADMINISTRASI" table.
# This is the 2nd table in the document. Row 3 = item "1" (tata ruang / RTRW),
# row 4 = item "2" (persetujuan awal). We insert a brand-new row before the
# current row 4 so it becomes the new item "2" (PIPPIB), pushing every
# following numbered row down by one position. Word clones the cell
# structure/formatting of the reference row automatically, so we only need
# to fill in the text afterwards and bump the "NO" column of every row that
# shifted down.

$d = $word.ActiveDocument
$t = $d.Tables(2)

$dollar = [char]36

# Insert a new, blank row immediately above the current row 4
# (the "persetujuan_awal" / item "2" row). Word duplicates the formatting
# of row 4 for the freshly inserted row.
$refRow = $t.Rows(4)
$t.Rows.Add($refRow) | Out-Null

# Populate the newly inserted row (now row 4) with the PIPPIB content.
$newRow = $t.Rows(4)
$newRow.Cells(1).Range.Text = "2"
$newRow.Cells(2).Range.Text = "Justifikasi/bukti kesesuaian lokasi rencana usaha dan/atau kegiatan dengan PIPPIB"
$newRow.Cells(3).Range.Text = $dollar + "{pippib_exist}"
$newRow.Cells(4).Range.Text = $dollar + "{pippib_not_exist}"
$newRow.Cells(5).Range.Text = $dollar + "{pippib_yes}"
$newRow.Cells(6).Range.Text = $dollar + "{pippib_no}"
$newRow.Cells(7).Range.Text = $dollar + "{pippib_ket}"

# Renumber every row below the new one: the old items "2".."10" (now sitting
# in rows 5..13) keep their own text/placeholders, but their displayed
# number must increase by one (2->3, 3->4, ... 10->11).
for ($i = $t.Rows.Count; $i -ge 5; $i--) {
    $row = $t.Rows($i)
    $numText = $row.Cells(1).Range.Text
    $numText = $numText -replace "[\x07\x0d\x0a]", ""
    if ($numText -match '^[0-9]+$') {
        $newNum = [int]$numText + 1
        $row.Cells(1).Range.Text = [string]$newNum
    }
}
